$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coefficients")

# --- New "MEANS" block (rows 79-83) ---
# Written first so the new shared strings ("MEANS", "commune", "SD")
# are minted in the same order as the target workbook.
$ws.Range("B79").Value = "MEANS"
$ws.Range("C79").Value = "Coeff"

$ws.Range("B80").Value = "commune"
$ws.Range("C80").Formula = "=AVERAGE(C74:D74)"
$ws.Range("D80").Formula = "=AVERAGE(E74:F74)"

$ws.Range("D79").Value = "SD"

$ws.Range("B81").Value = "com_year"
$ws.Range("C81").Formula = "=AVERAGE(C75:D75)"
$ws.Range("D81").Formula = "=AVERAGE(E75:F75)"

$ws.Range("B82").Value = "Province"
$ws.Range("C82").Formula = "=AVERAGE(C76:D76)"
$ws.Range("D82").Formula = "=AVERAGE(E76:F76)"

$ws.Range("B83").Value = "prov_year"
$ws.Range("C83").Formula = "=AVERAGE(C77:D77)"
$ws.Range("D83").Formula = "=AVERAGE(E77:F77)"

# --- New mini "var/SD" summary table (H45:J47), mirroring H39:K43 ---
$ws.Range("I45").Value = "var"
$ws.Range("J45").Value = "SD"

$ws.Range("H46").Value = "prov"
$ws.Range("I46").Formula = "=AVERAGE(J40,J42)"
$ws.Range("J46").Formula = "=AVERAGE(K40,K42)"

$ws.Range("H47").Value = "year"
$ws.Range("I47").Formula = "=AVERAGE(J41,J43)"
$ws.Range("J47").Formula = "=AVERAGE(K41,K43)"

# --- Restore view state: active sheet, scroll position, selection ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 14 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("M31").Select()
